$d = $word.ActiveDocument

$replacements = @(
    @{old = "2026-02-08 Sunday"; new = "2026-02-09 Monday"},
    @{old = "472×6="; new = "269×9="},
    @{old = "678×9="; new = "374×8="},
    @{old = "817×5="; new = "820×6="},
    @{old = "860×7="; new = "600×3="},
    @{old = "366×6="; new = "640×8="},
    @{old = "913×7="; new = "165×8="},
    @{old = "256×8="; new = "875×2="},
    @{old = "717×9="; new = "632×5="},
    @{old = "584×4="; new = "167×6="},
    @{old = "226×9="; new = "991×8="},
    @{old = "728×4="; new = "927×5="},
    @{old = "648×5="; new = "536×5="},
    @{old = "744×3="; new = "842×7="},
    @{old = "612×7="; new = "322×4="},
    @{old = "947×3="; new = "252×4="},
    @{old = "742×6="; new = "599×5="},
    @{old = "275×9="; new = "136×3="},
    @{old = "234×8="; new = "726×3="},
    @{old = "727×4="; new = "249×7="},
    @{old = "666×4="; new = "611×8="},
    @{old = "283×7="; new = "865×7="},
    @{old = "195×2="; new = "571×9="},
    @{old = "879×5="; new = "550×2="},
    @{old = "213×5="; new = "332×8="},
    @{old = "299×7="; new = "968×5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
